$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a "last changed" date serial that was bumped
# from 45206 (2023-10-07) to 45208 (2023-10-09) for every data row.
$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45206) {
        $cell.Value2 = 45208
    }
}
